$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SGD")
$ws.Range("B2").Value = 0.5600000000000001
$ws.Range("C2").Value = 0.625
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.631578947368421
$ws.Range("F2").Value = 0.7283333333333334
$ws.Range("G2").Value = 0.7293421052631579
$ws.Range("C3").Value = 0.4545454545454545
$ws.Range("D3").Value = 0.3846153846153846
$ws.Range("E3").Value = 0.631578947368421
$ws.Range("F3").Value = 0.6130536130536131
$ws.Range("G3").Value = 0.631578947368421
$ws.Range("B4").Value = 0.717948717948718
$ws.Range("C4").Value = 0.5263157894736842
$ws.Range("D4").Value = 0.5555555555555556
$ws.Range("E4").Value = 0.631578947368421
$ws.Range("F4").Value = 0.5999400209926525
$ws.Range("G4").Value = 0.6069204725714421
$ws.Range("B5").Value = 14
$ws.Range("D5").Value = 13
$ws.Range("E5").Value = 0.631578947368421

$ws = $wb.Worksheets.Item("LinearSVC")
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("B5").Value = 14
$ws.Range("D5").Value = 13
$ws.Range("E5").Value = 1

$ws = $wb.Worksheets.Item("MLP Neural Network")
$ws.Range("D2").Value = 0.9285714285714286
$ws.Range("F2").Value = 0.9761904761904763
$ws.Range("G2").Value = 0.9755639097744361
$ws.Range("D4").Value = 0.962962962962963
$ws.Range("F4").Value = 0.9717813051146384
$ws.Range("G4").Value = 0.9735449735449735
$ws.Range("B5").Value = 14
$ws.Range("D5").Value = 13

$ws = $wb.Worksheets.Item("Gaussian Process")
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("B5").Value = 14
$ws.Range("D5").Value = 13
$ws.Range("E5").Value = 1
